$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Text replacements (order-independent; unique search strings)
# ------------------------------------------------------------------

# "Philadelphia Crime/Weather Data:" -> "...and Philadelphia Parking Violation:"
$d.Content.Find.Execute(
    "Philadelphia Crime/Weather Data:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Philadelphia Crime/Weather Data and Philadelphia Parking Violation:",
    2
) | Out-Null

# "Data cleansing would include..." -> new "From last term..." paragraph text
$d.Content.Find.Execute(
    "Data cleansing would include removing of rows with null property values to have zip-codes in the dataset which matched the crime data set.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "From last term, these Philadelphia Crime/Weather Data were cleansed and integrated. We will integrate Philadelphia Parking Violation based on location features such as lat and lon. From this we will be able to create zip codes for crime predictive analysis. ",
    2
) | Out-Null

# Machine learning paragraph rewrite
$d.Content.Find.Execute(
    "The machine learning part of the project will be made up of three components: data transformation/normalization, regression, and classification. For classification, we will do our best to implement all of the following model types: Logistic Regression, Naive Bayes, SVM, Decision Tree Classifier, Randomforest, Neural Networks, LDA, QDA, KNN, NN. The next part is regression, which will consist of all of the following: Linear Regression, add regularizer to linear regression (Ridge and Lasso Regression), Decision Tree Regression, Support Vector Regression, Xgboost, Neural Networks (NN), partial least square regression, In addition, for data transformation and normalization, we will perform at least one of the following: Dimensionality reduction (PCA as an example), Hyper-parameter tuning (for example gamma and C for SVM), and Feature Selection. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The machine learning part of the project will be made up of three components: data transformation/normalization, data integration, and classification. For classification, we will do our best to implement all of the following model types: Logistic Regression, Decision Tree Classifier, Randomforest, Extra Trees, and KNN. In addition, for data transformation and normalization, we hope to perform at least one of the following: Dimensionality reduction (PCA as an example), Hyper-parameter tuning (for example gamma and C for SVM), and Feature Selection. ",
    2
) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the three paragraphs that used to sit between the
#    "Philadelphia Crime/Weather Data..." heading and the bullet that
#    now reads "From last term, these Philadelphia Crime/Weather
#    Data...":
#      - "From last term, these two datasets were cleansed..." bullet
#      - the blank paragraph right after it
#      - "Philadelphia Housing Data:" paragraph
#    Locate indices fresh (paragraph count is unaffected by the text
#    replacements above) and delete from the highest index down so
#    earlier indices stay valid.
# ------------------------------------------------------------------

$idxOldBullet = -1
$idxHousing = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("From last term, these two datasets were cleansed")) {
        $idxOldBullet = $i
    }
    if ($t.Contains("Philadelphia Housing Data:")) {
        $idxHousing = $i
    }
}

if ($idxHousing -gt 0) {
    $d.Paragraphs.Item($idxHousing).Range.Delete() | Out-Null
}
if ($idxOldBullet -gt 0) {
    # the blank paragraph is immediately after the old bullet paragraph
    $d.Paragraphs.Item($idxOldBullet + 1).Range.Delete() | Out-Null
    $d.Paragraphs.Item($idxOldBullet).Range.Delete() | Out-Null
}

# ------------------------------------------------------------------
# 3) The paragraph that used to hold only a tab character loses the
#    tab and gains ind(left=720, firstLine=0).
# ------------------------------------------------------------------

$idxTab = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "`t") {
        $idxTab = $i
        break
    }
}

if ($idxTab -gt 0) {
    $p = $d.Paragraphs.Item($idxTab)
    $p.LeftIndent = 36
    $p.FirstLineIndent = 0
    $rng = $d.Range($p.Range.Start, $p.Range.Start + 1)
    $rng.Delete() | Out-Null
}

Write-Output "done"
